# Apply updated odds values to row 3 of the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value  = 3.3
$ws.Range("K3").Value  = 2.05
$ws.Range("Q3").Value  = 2.2
$ws.Range("R3").Value  = 1.65
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6
$ws.Range("AH3").Value = 17
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 12
$ws.Range("AS3").Value = 201
$ws.Range("AU3").Value = 8.5
$ws.Range("BB3").Value = 251
